$d = $word.ActiveDocument
$wns = 'xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"'

function Find-ParagraphContaining($doc, [string]$needle) {
    for ($i = 1; $i -le $doc.Paragraphs.Count; $i++) {
        $p = $doc.Paragraphs($i)
        if ($p.Range.Text -like "*$needle*") {
            return $p
        }
    }
    return $null
}

# ---------------------------------------------------------------------------
# Change 1: Pre-conditions bullet "Crisis Management main page is opened."
#           -> "Crisis Board page is opened."
# ---------------------------------------------------------------------------
$crisisPara = Find-ParagraphContaining $d "Crisis Management main page is opened."
if ($null -eq $crisisPara) {
    throw "Could not locate the 'Crisis Management main page is opened.' paragraph"
}

$crisisXml = '<w:p ' + $wns + ' w:rsidR="006660D5" w:rsidRPr="001E679A" w:rsidRDefault="006660D5" w:rsidP="006660D5">' +
    '<w:pPr>' +
        '<w:pStyle w:val="NormalWeb"/>' +
        '<w:numPr><w:ilvl w:val="0"/><w:numId w:val="3"/></w:numPr>' +
        '<w:spacing w:before="0" w:after="0"/>' +
        '<w:rPr><w:color w:val="000000"/><w:sz w:val="20"/><w:szCs w:val="20"/></w:rPr>' +
    '</w:pPr>' +
    '<w:r><w:rPr><w:color w:val="000000"/><w:sz w:val="20"/><w:szCs w:val="20"/></w:rPr>' +
        '<w:t xml:space="preserve">Crisis </w:t></w:r>' +
    '<w:r><w:rPr><w:color w:val="000000"/><w:sz w:val="20"/><w:szCs w:val="20"/></w:rPr>' +
        '<w:t>Board</w:t></w:r>' +
    '<w:r><w:rPr><w:color w:val="000000"/><w:sz w:val="20"/><w:szCs w:val="20"/></w:rPr>' +
        '<w:t xml:space="preserve"> page is opened.</w:t></w:r>' +
    '</w:p>'

$crisisPara.Range.InsertXML($crisisXml)

# ---------------------------------------------------------------------------
# Change 2: Main Path step 1 "Manager selects the Incident from incidents
#           list and presses 'Close Incident'" is reworded and four new
#           steps are added describing the new navigation flow.
# ---------------------------------------------------------------------------
$stepPara = Find-ParagraphContaining $d "selects the Incident from incidents list"
if ($null -eq $stepPara) {
    throw "Could not locate the 'Manager selects the Incident...' paragraph"
}

$rpr = '<w:rPr><w:color w:val="000000"/><w:sz w:val="20"/><w:szCs w:val="20"/></w:rPr>'
$pPr = '<w:pPr><w:pStyle w:val="NormalWeb"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="5"/></w:numPr>' +
       '<w:spacing w:before="0" w:after="0"/>' + $rpr + '</w:pPr>'

$stepXml =
    '<w:p ' + $wns + ' w:rsidR="006660D5" w:rsidRDefault="006660D5" w:rsidP="006660D5">' + $pPr +
        '<w:r>' + $rpr + '<w:t xml:space="preserve">Manager selects the </w:t></w:r>' +
        '<w:r>' + $rpr + '<w:t>List all i</w:t></w:r>' +
        '<w:r>' + $rpr + '<w:t>ncident</w:t></w:r>' +
        '<w:r>' + $rpr + '<w:t>s</w:t></w:r>' +
        '<w:r>' + $rpr + '<w:t xml:space="preserve"> from </w:t></w:r>' +
        '<w:r>' + $rpr + '<w:t>I</w:t></w:r>' +
        '<w:r>' + $rpr + '<w:t>ncident</w:t></w:r>' +
        '<w:r>' + $rpr + '<w:t xml:space="preserve"> </w:t></w:r>' +
        '<w:r>' + $rpr + '<w:t>menu panel.</w:t></w:r>' +
    '</w:p>' +
    '<w:p ' + $wns + '>' + $pPr +
        '<w:r>' + $rpr + '<w:t>System shows the Incident List page with the list of the incidents.</w:t></w:r>' +
    '</w:p>' +
    '<w:p ' + $wns + '>' + $pPr +
        '<w:r>' + $rpr + '<w:t>Manager selects the incident from the list.</w:t></w:r>' +
    '</w:p>' +
    '<w:p ' + $wns + '>' + $pPr +
        '<w:r>' + $rpr + '<w:t>System shows Edit Incident page with the information of the incident.</w:t></w:r>' +
    '</w:p>' +
    '<w:p ' + $wns + '>' + $pPr +
        '<w:r>' + $rpr + '<w:t>Manager selects Close button.</w:t></w:r>' +
    '</w:p>'

$stepPara.Range.InsertXML($stepXml)
